# GLES_VariableInfo.xlsx update:
# Insert a new "item_num" column just before the existing "comment" column
# on the Measures, ID, Dems, Dates and NewVars sheets; populate it with 1
# on the Measures sheet (the only sheet that already had row data past
# column L). Update the ID!_FilterDatabase defined name to cover the new
# column, and move the active-sheet/selection focus onto NewVars.

$wb = $excel.ActiveWorkbook

# --- Measures sheet: comment was column Q (17); new item_num goes there,
# comment shifts to R (18); fill item_num = 1 for the 8 data rows.
$ws = $wb.Worksheets.Item("Measures")
$ws.Columns.Item(17).Insert()
$ws.Range("Q1").Value = "item_num"
$ws.Range("Q2").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("Q10").Select()

# --- ID sheet: comment was column Q (17); new item_num goes there, comment
# shifts to R (18). No data in the new column besides the header.
$ws = $wb.Worksheets.Item("ID")
$ws.Columns.Item(17).Insert()
$ws.Range("Q1").Value = "item_num"
$ws.Range("Q2").Select()

# --- Dems sheet: same pattern as ID.
$ws = $wb.Worksheets.Item("Dems")
$ws.Columns.Item(17).Insert()
$ws.Range("Q1").Value = "item_num"
$ws.Range("Q2").Select()

# --- Dates sheet: same pattern as ID.
$ws = $wb.Worksheets.Item("Dates")
$ws.Columns.Item(17).Insert()
$ws.Range("Q1").Value = "item_num"
$ws.Range("Q2").Select()

# --- NewVars sheet: comment was column P (16); new item_num goes there,
# comment shifts to Q (17). This sheet becomes the active tab/sheet.
$ws = $wb.Worksheets.Item("NewVars")
$ws.Columns.Item(16).Insert()
$ws.Range("P1").Value = "item_num"
$ws.Activate()
$ws.Range("P2").Select()

# --- Update the ID sheet's hidden AutoFilter defined name so it spans the
# new column (was $A$1:$Q$1, now $A$1:$R$1).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ID!_FilterDatabase") {
        $n.RefersTo = "=ID!`$A`$1:`$R`$1"
    }
}
